$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds vaccine-combination labels in ALL CAPS (e.g. MODERNA_JANSSEN,
# PFIZER_JANSSEN, JANSSEN, PFIZER_MODERNA, MODERNA, PFIZER) for every data row
# except the header (row 1) and the "Shared" rows (rows 2-3), which are left
# untouched. Convert each underscore-separated word to Title Case, e.g.
# MODERNA_JANSSEN -> Moderna_Janssen, JANSSEN -> Janssen.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 4; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2

    $parts = $val.ToString().Split("_")
    $newParts = @()
    foreach ($p in $parts) {
        if ($p.Length -gt 0) {
            $newParts += ($p.Substring(0,1).ToUpper() + $p.Substring(1).ToLower())
        } else {
            $newParts += $p
        }
    }
    $newVal = [string]::Join("_", $newParts)

    $cell.Value = $newVal
}
